{"js": "const replacements = [\n  [\"2026-01-04 Sunday\", \"2026-01-05 Monday\"],\n  [\"823\u00d74=3292\", \"746\u00d74=2984\"],\n  [\"627\u00d78=5016\", \"653\u00d72=1306\"],\n  [\"796\u00d75=3980\", \"177\u00d77=1239\"],\n  [\"285\u00d73=855\", \"624\u00d79=5616\"],\n  [\"998\u00d77=6986\", \"806\u00d73=2418\"],\n  [\"212\u00d73=636\", \"200\u00d77=1400\"],\n  [\"972\u00d75=4860\", \"609\u00d73=1827\"],\n  [\"855\u00d78=6840\", \"885\u00d75=4425\"],\n  [\"911\u00d75=4555\", \"927\u00d72=1854\"],\n  [\"271\u00d78=2168\", \"355\u00d76=2130\"],\n  [\"733\u00d77=5131\", \"559\u00d78=4472\"],\n  [\"339\u00d78=2712\", \"937\u00d78=7496\"],\n  [\"560\u00d78=4480\", \"209\u00d79=1881\"],\n  [\"179\u00d73=537\", \"618\u00d76=3708\"],\n  [\"725\u00d76=4350\", \"795\u00d72=1590\"],\n  [\"462\u00d74=1848\", \"578\u00d72=1156\"],\n  [\"323\u00d74=1292\", \"924\u00d76=5544\"],\n  [\"361\u00d79=3249\", \"180\u00d73=540\"],\n  [\"417\u00d75=2085\", \"878\u00d77=6146\"],\n  [\"363\u00d75=1815\", \"262\u00d74=1048\"],\n  [\"852\u00d78=6816\", \"988\u00d72=1976\"],\n  [\"144\u00d77=1008\", \"483\u00d74=1932\"],\n  [\"819\u00d76=4914\", \"224\u00d73=672\"],\n  [\"912\u00d78=7296\", \"365\u00d73=1095\"],\n  [\"243\u00d72=486\", \"331\u00d78=2648\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-04 Sunday\", \"2026-01-05 Monday\"),\n    @(\"823\u00d74=3292\", \"746\u00d74=2984\"),\n    @(\"627\u00d78=5016\", \"653\u00d72=1306\"),\n    @(\"796\u00d75=3980\", \"177\u00d77=1239\"),\n    @(\"285\u00d73=855\", \"624\u00d79=5616\"),\n    @(\"998\u00d77=6986\", \"806\u00d73=2418\"),\n    @(\"212\u00d73=636\", \"200\u00d77=1400\"),\n    @(\"972\u00d75=4860\", \"609\u00d73=1827\"),\n    @(\"855\u00d78=6840\", \"885\u00d75=4425\"),\n    @(\"911\u00d75=4555\", \"927\u00d72=1854\"),\n    @(\"271\u00d78=2168\", \"355\u00d76=2130\"),\n    @(\"733\u00d77=5131\", \"559\u00d78=4472\"),\n    @(\"339\u00d78=2712\", \"937\u00d78=7496\"),\n    @(\"560\u00d78=4480\", \"209\u00d79=1881\"),\n    @(\"179\u00d73=537\", \"618\u00d76=3708\"),\n    @(\"725\u00d76=4350\", \"795\u00d72=1590\"),\n    @(\"462\u00d74=1848\", \"578\u00d72=1156\"),\n    @(\"323\u00d74=1292\", \"924\u00d76=5544\"),\n    @(\"361\u00d79=3249\", \"180\u00d73=540\"),\n    @(\"417\u00d75=2085\", \"878\u00d77=6146\"),\n    @(\"363\u00d75=1815\", \"262\u00d74=1048\"),\n    @(\"852\u00d78=6816\", \"988\u00d72=1976\"),\n    @(\"144\u00d77=1008\", \"483\u00d74=1932\"),\n    @(\"819\u00d76=4914\", \"224\u00d73=672\"),\n    @(\"912\u00d78=7296\", \"365\u00d73=1095\"),\n    @(\"243\u00d72=486\", \"331\u00d78=2648\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
